# Append the new monthly data row (2024-10-01) to the bottom of the table,
# mirroring the existing rows' layout (Date, NEIG, Food, Energy, Rent, Core Services).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 58
$newRow = 59

# Copy the formatting (incl. the date number format / font on column A) from
# the last existing row down onto the new row before filling in values, so we
# reuse the workbook's existing style instead of minting a new one.
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy() | Out-Null
$ws.Range("A" + $newRow + ":F" + $newRow).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($newRow, 1).Value = 45566   # 2024-10-01
$ws.Cells.Item($newRow, 2).Value = -0.458  # NEIG
$ws.Cells.Item($newRow, 3).Value = 0.462   # Food
$ws.Cells.Item($newRow, 4).Value = -1.428  # Energy
$ws.Cells.Item($newRow, 5).Value = 0.369   # Rent
$ws.Cells.Item($newRow, 6).Value = 1.211   # Core Services
